$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    # Force text interpretation for numeric-looking strings so that
    # values like "376.49" are preserved exactly as text, not coerced to a number.
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range('D2').Value = '51.098.88'
$ws.Range('E2').Value = '  -1.10%  '
$ws.Range('D3').Value = '2.939.34'
$ws.Range('E3').Value = '  -1.90%  '
$ws.Range('E4').Value = '  -0.10%  '
Set-TextValue $ws 'D5' '376.49'
$ws.Range('E5').Value = '  -1.29%  '
$ws.Range('E6').Value = '  -3.57%  '
$ws.Range('E7').Value = '  -2.00%  '
$ws.Range('E8').Value = '  -0.09%  '
Set-TextValue $ws 'D9' '0.584'
$ws.Range('E9').Value = '  -2.54%  '
Set-TextValue $ws 'D10' '36.70'
$ws.Range('E10').Value = '  -2.41%  '
$ws.Range('E11').Value = '  -1.01%  '
Set-TextValue $ws 'D12' '0.0835'
$ws.Range('E12').Value = '  -1.13%  '
$ws.Range('D13').Value = '3.401.26'
$ws.Range('E13').Value = '  -2.15%  '
Set-TextValue $ws 'D14' '17.95'
$ws.Range('E14').Value = '  -4.13%  '
$ws.Range('E15').Value = '  -2.44%  '
$ws.Range('D16').Value = '2.938.41'
$ws.Range('E16').Value = '  -1.85%  '
Set-TextValue $ws 'D17' '0.973'
$ws.Range('E17').Value = '  +0.13%  '
$ws.Range('D18').Value = '51.038.14'
$ws.Range('E18').Value = '  -1.36%  '
Set-TextValue $ws 'D19' '3.15'
$ws.Range('E19').Value = '  -6.79%  '
Set-TextValue $ws 'D20' '7.13'
$ws.Range('E20').Value = '  -4.04%  '
Set-TextValue $ws 'D21' '12.54'
$ws.Range('E21').Value = '  -4.20%  '
$ws.Range('D22').Value = '0.0₃0954'
$ws.Range('E22').Value = '  -0.65%  '
Set-TextValue $ws 'D23' '263.04'
$ws.Range('E23').Value = '  -0.32%  '
Set-TextValue $ws 'D24' '68.19'
$ws.Range('E24').Value = '  -1.72%  '
Set-TextValue $ws 'D25' '2.87'
$ws.Range('E25').Value = '  +2.50%  '
Set-TextValue $ws 'D26' '8.14'
$ws.Range('E26').Value = '  +9.01%  '
Set-TextValue $ws 'D27' '7.82'
$ws.Range('E27').Value = '  +8.13%  '
$ws.Range('B28').Value = 'Hedera'
$ws.Range('C28').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws 'D28' '0.113'
$ws.Range('E28').Value = '  +5.32%  '
$ws.Range('B29').Value = 'Dai'
$ws.Range('C29').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws 'D29' '1.00'
$ws.Range('E29').Value = '  +0.00%  '
Set-TextValue $ws 'D30' '0.167'
$ws.Range('E30').Value = '  -2.70%  '
Set-TextValue $ws 'D31' '25.67'
$ws.Range('E31').Value = '  -1.87%  '
$ws.Range('E32').Value = '  -0.79%  '
Set-TextValue $ws 'D33' '34.07'
$ws.Range('E33').Value = '  -2.25%  '
Set-TextValue $ws 'D34' '50.73'
$ws.Range('E34').Value = '  -1.34%  '
Set-TextValue $ws 'D35' '0.0455'
$ws.Range('E35').Value = '  -0.65%  '
$ws.Range('E36').Value = '  -3.43%  '
$ws.Range('E37').Value = '  -0.11%  '
$ws.Range('E38').Value = '  -4.99%  '
$ws.Range('E39').Value = '  -2.57%  '
Set-TextValue $ws 'D40' '16.47'
$ws.Range('E40').Value = '  -6.13%  '
$ws.Range('E41').Value = '  -1.68%  '
$ws.Range('E42').Value = '  -4.14%  '
Set-TextValue $ws 'D43' '121.53'
$ws.Range('E43').Value = '  -2.09%  '
Set-TextValue $ws 'D44' '21.11'
$ws.Range('E44').Value = '  -5.45%  '
$ws.Range('E45').Value = '  -1.93%  '
Set-TextValue $ws 'D46' '0.270'
$ws.Range('E46').Value = '  -1.87%  '
$ws.Range('E47').Value = '  +0.43%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws 'D48' '3.22'
$ws.Range('E48').Value = '  -1.65%  '
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').Value = '2.001.93'
$ws.Range('E49').Value = '  -2.89%  '
Set-TextValue $ws 'D50' '0.0345'
$ws.Range('E50').Value = '  -3.02%  '
$ws.Range('E51').Value = '  -3.48%  '

Write-Output "Applied 92 cell updates"
